{"js": "// Update the \"classe Utente\" sentence in the domain-model description:\n// the \"(opzionale)\" marker moves from \"telefono\" to \"data di nascita\".\n//\n//   before: \"... telefono (opzionale), data di nascita ed un'immagine di profilo (opzionale).\"\n//   after : \"... telefono, data di nascita (opzionale) ed un'immagine di profilo (opzionale).\"\nconst searchText = \"telefono (opzionale), data di nascita ed un\\u2019immagine di profilo\";\nconst replacement = \"telefono, data di nascita (opzionale) ed un\\u2019immagine di profilo\";\n\nconst results = context.document.body.search(searchText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target sentence fragment not found in document body.\");\n}\n\nresults.items[0].insertText(replacement, \"Replace\");\nawait context.sync();\n", "ps1": "# Update the \"classe Utente\" sentence in the domain-model description:\n# the \"(opzionale)\" marker moves from \"telefono\" to \"data di nascita\".\n#\n#   before: \"... telefono (opzionale), data di nascita ed un'immagine di profilo (opzionale).\"\n#   after : \"... telefono, data di nascita (opzionale) ed un'immagine di profilo (opzionale).\"\n$d = $word.ActiveDocument\n\n$rightQuote = [char]0x2019\n$findText = \"telefono (opzionale), data di nascita ed un\" + $rightQuote + \"immagine di profilo\"\n$replaceText = \"telefono, data di nascita (opzionale) ed un\" + $rightQuote + \"immagine di profilo\"\n\n$range = $d.Content\n# wdFindContinue=1 (Wrap), wdReplaceOne=1 (Replace) - only the single, unique match is changed.\n$found = $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1)\n\nif (-not $found) {\n    throw \"Target sentence fragment not found in document.\"\n}\n"}
